# ---------------------------------------------------------------------------
# Applies the commit "Added a few more slots" to the Cai Shen Dao review doc:
#   1. Inserts a new "Meta description: ..." paragraph right after the
#      Heading1 title paragraph.
#   2. Removes the duplicate bold title paragraph that was mistakenly left
#      near the end of the document.
#   3. Replaces the text of the final (italic) paragraph - formerly holding
#      a copy of the meta description - with the image-generation prompt.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
# Paragraph.Range.Text always carries a trailing paragraph-mark character.
$cr = [char]13

# --- 1. Insert the "Meta description" paragraph after the title -----------
$titlePara = $d.Paragraphs.Item(1)
if ($titlePara.Range.Text -ne ("Play Cai Shen Dao slot free: Review and analysis" + $cr)) {
    throw "Unexpected paragraph 1 text: $($titlePara.Range.Text)"
}
$insertion = $d.Range($titlePara.Range.End, $titlePara.Range.End)

# Build the new paragraph's raw OOXML so we can reproduce the document's
# existing convention of a leading empty run, a bold "Meta description"
# run, and a plain run carrying the rest of the sentence.
$metaParagraphXml = "<w:p $wNs>" +
    "<w:r/>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>" +
    "<w:r><w:t>: Explore Cai Shen Dao slot game with our review. Play for free and discover its impressive graphics, thematic symbols, medium volatility, and free spin feature.</w:t></w:r>" +
    "</w:p>"
# A trailing empty paragraph is required so InsertXML treats the paragraph
# above as a genuine standalone block instead of merging its runs into the
# following paragraph; it is deleted again immediately afterwards.
$trailingMarkerXml = "<w:p $wNs/>"

[void]$insertion.InsertXML($metaParagraphXml + $trailingMarkerXml)

$markerPara = $d.Paragraphs.Item(3)
if ($markerPara.Range.Text -ne ("" + $cr)) {
    throw "Unexpected marker paragraph text: $($markerPara.Range.Text)"
}
$markerPara.Range.Delete()

# --- 2. Remove the duplicate bold title paragraph near the end ------------
$paraCount = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($paraCount - 1)
if ($dupTitlePara.Range.Text -ne ("Play Cai Shen Dao slot free: Review and analysis" + $cr)) {
    throw "Unexpected duplicate-title paragraph text: $($dupTitlePara.Range.Text)"
}
$dupTitlePara.Range.Delete()

# --- 3. Replace the final paragraph's text with the image prompt ----------
$paraCount = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($paraCount)
$oldMetaText = "Explore Cai Shen Dao slot game with our review. Play for free and discover its impressive graphics, thematic symbols, medium volatility, and free spin feature."
if ($finalPara.Range.Text -ne ($oldMetaText + $cr)) {
    throw "Unexpected final paragraph text: $($finalPara.Range.Text)"
}
$promptText = "Prompt: Create a fun and engaging feature image for `"Cai Shen Dao`" that fits the description of a happy Maya warrior with glasses. The image should be in cartoon style and showcase the Chinese culture theme of the game, with symbols such as the carp, fan, and amulets. The image should also have a touch of humor and a playful vibe to appeal to the game's target audience."

$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End)
$finalRange.Text = $promptText
